# Update "想去人数" (want-to-go count) values in column F
# for worksheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 840
$ws1.Range("F11").Value = 2090
$ws1.Range("F13").Value = 1521
$ws1.Range("F14").Value = 2713
$ws1.Range("F16").Value = 4040
$ws1.Range("F18").Value = 158
$ws1.Range("F22").Value = 30
$ws1.Range("F26").Value = 3790
$ws1.Range("F28").Value = 3290
$ws1.Range("F29").Value = 1102
$ws1.Range("F34").Value = 260
$ws1.Range("F35").Value = 383
$ws1.Range("F36").Value = 245

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 840
$ws4.Range("F13").Value = 2090
$ws4.Range("F15").Value = 1521
$ws4.Range("F17").Value = 2713
$ws4.Range("F19").Value = 4040
$ws4.Range("F21").Value = 158
$ws4.Range("F25").Value = 30
$ws4.Range("F30").Value = 3790
$ws4.Range("F32").Value = 3290
$ws4.Range("F33").Value = 1102
$ws4.Range("F38").Value = 260
$ws4.Range("F39").Value = 383
$ws4.Range("F40").Value = 245
